$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Table on slide 6: change the table style (tableStyleId) from
#    {E400C113-A2C2-4DA8-B2BF-A90CE7AC23B6} to
#    {1F3C3EAB-F919-4852-AA43-753D6A7E6D31}
# ------------------------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{1F3C3EAB-F919-4852-AA43-753D6A7E6D31}")
    }
}

# ------------------------------------------------------------------
# 2) Theme colours: swap the "Integral" colour scheme (currently on
#    the deck's master theme) for the stock "Office Theme" colour
#    scheme values.
# ------------------------------------------------------------------
function Set-ThemeColor($scheme, $index, $hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Item($index).RGB = $bb * 65536 + $gg * 256 + $rr
}

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
Set-ThemeColor $masterScheme 1  "000000"
Set-ThemeColor $masterScheme 2  "FFFFFF"
Set-ThemeColor $masterScheme 3  "44546A"
Set-ThemeColor $masterScheme 4  "E7E6E6"
Set-ThemeColor $masterScheme 5  "5B9BD5"
Set-ThemeColor $masterScheme 6  "ED7D31"
Set-ThemeColor $masterScheme 7  "A5A5A5"
Set-ThemeColor $masterScheme 8  "FFC000"
Set-ThemeColor $masterScheme 9  "4472C4"
Set-ThemeColor $masterScheme 10 "70AD47"
Set-ThemeColor $masterScheme 11 "0563C1"
Set-ThemeColor $masterScheme 12 "954F72"
